$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (14), shifting existing N:P data to O:Q
$ws.Range("N1").EntireColumn.Insert()

# The newly inserted column inherits the width of the column to its left (M)
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Update selection to reflect where the cursor ended up after the edit
$ws.Range("R6").Select()
